$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current AS (column 45, "Ly do xoa bo"),
# right after AR ("Ma tra cuu"). This pushes the remaining trailing
# headers (Ly do xoa bo, Thong tin hoa don lien quan, Ngay lap, Nguoi lap)
# one column to the right.
$ws.Columns.Item(45).Insert()

# Give the brand-new column (now AS) its header text and matching width.
$ws.Range("AS4").Value = "Mã CQT cấp"
$ws.Columns.Item(45).ColumnWidth = 37.125

# Rename header AQ4 from "Trang thai phat hanh" to "Trang thai quy trinh"
$ws.Range("AQ4").Value = "Trạng thái quy trình"

# Adjust the sheet view: hide gridlines, scroll to show the new trailing
# columns, and move the active selection to A5.
$ws.Application.ActiveWindow.DisplayGridlines = $false
$ws.Application.Goto($ws.Range("AO1"))
$ws.Range("A5").Select()
